$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. "Advantages:" paragraph -- shorten the sentence and relocate the
#    (hidden) "_GoBack" bookmark to sit right after the new wording.
#    A throwaway marker ("ZZZMARK") is appended so the bookmark can be
#    anchored at a safe (non-paragraph-boundary) offset, then the
#    marker text is deleted, leaving the bookmark collapsed exactly
#    where we need it.
# ---------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Execute("are being set up in the constructor. ", $false, $false, $false, $false, $false, $true, 1, $false, "are being initialised in the constructor having access to those private variablesZZZMARK", 2)

$markerRange = $d.Content
$markerRange.Find.Execute("ZZZMARK")
$d.Bookmarks.Add("_GoBack", $d.Range($markerRange.Start, $markerRange.Start))
$markerRange.Delete()

# ---------------------------------------------------------------------
# 2. "Disadvantages:" -> "Problems:"
# ---------------------------------------------------------------------
$find2 = $d.Content.Find
$find2.ClearFormatting()
$find2.Execute("Disadvantages:", $false, $false, $false, $false, $false, $true, 1, $false, "Problems:", 2)

# ---------------------------------------------------------------------
# 3. Replace the (previously yellow-highlighted) "The constructor for
#    the random number generator will be called twice? " paragraph
#    with the new "However, the rng_ constructor ..." paragraph, and
#    strip the old highlighting.
# ---------------------------------------------------------------------
$oldParaText = "The constructor for the random number generator will be called twice? "
$newParaText = "However, the rng_ constructor won" + [char]0x2019 + "t have been run straight away. This will only happen when the symbol, x, y, and the mouse pointer have been initialised. This means that x and y value for the snake will not be random. Only time it will be random is when the set up function for the game is run due to the position at random function in that game function."

$find3 = $d.Content.Find
$find3.ClearFormatting()
$find3.Execute($oldParaText, $false, $false, $false, $false, $false, $true, 1, $false, $newParaText, 2)

$find4 = $d.Content
$find4.Find.ClearFormatting()
$find4.Find.Execute($newParaText)
$find4.HighlightColorIndex = 0

# ---------------------------------------------------------------------
# 4. Add the new closing sentence as its own paragraph straight after.
# ---------------------------------------------------------------------
$find5 = $d.Content
$find5.Find.ClearFormatting()
$find5.Find.Execute($newParaText)
$targetStart = $find5.Start

$howeverIndex = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -le $targetStart -and $targetStart -lt $p.Range.End) {
        $howeverIndex = $i
        break
    }
}

$find5.Collapse(0)
$find5.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($howeverIndex + 1)
$newPara.Range.InsertAfter("In other words, this will run the same function twice.")

Write-Output "done"
